# Rename header labels on the existing sheets.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$ws2 = $wb.Worksheets.Item(2)   # "Monthly Trend"

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after "Monthly Trend" (at the end).
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "PO Forecast"

# Reuse the header style (bold / bordered / centered) from the existing
# sheets so no new cellXfs entries get minted.
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# Reuse the date-column style (custom date/time number format) too.
$ws1.Range("A2").Copy()
$ws3.Range("A2:A28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

$rows = @(
    (@(45228.99999999999, 34, 4.740202526031666, 62.69479649522214)),
    (@(45270.99999999999, 33, 4.648621315659796, 63.0219760321971)),
    (@(45298.99999999999, 32, 4.424795512778714, 61.9263218202944)),
    (@(45305.99999999999, 32, 2.144176156808419, 60.56817044889288)),
    (@(45319.99999999999, 32, 2.250976812151276, 60.06018138909931)),
    (@(45333.99999999999, 32, 2.164509288752438, 61.6335683259378)),
    (@(45347.99999999999, 31, 2.953707760037067, 62.38171153366238)),
    (@(45361.99999999999, 31, 3.014459062605261, 59.05237286595398)),
    (@(45375.99999999999, 31, 2.593501061854954, 61.30447249403041)),
    (@(45382.99999999999, 31, 1.145488427352993, 60.51541356589137)),
    (@(45515.99999999999, 28, -2.165118148821175, 56.78798256750386)),
    (@(45522.99999999999, 28, 0.01644262266500549, 55.65694832535821)),
    (@(45529.99999999999, 28, -1.408540016946355, 57.44917852663067)),
    (@(45536.99999999999, 28, -0.4790627504794984, 54.69911785413857)),
    (@(45543.99999999999, 28, -0.5632406795680566, 58.20075606878463)),
    (@(45550.99999999999, 28, -2.663712925622445, 56.41453495178853)),
    (@(45557.99999999999, 27, -2.910748427598009, 54.92060906115279)),
    (@(45571.99999999999, 27, -1.202118811381786, 56.42677702447465)),
    (@(45599.99999999999, 27, -2.151412530723036, 52.86963612432499)),
    (@(45606.99999999999, 27, -2.706158130545303, 53.69364707632938)),
    (@(45613.99999999999, 26, -1.709061148657957, 54.08017869260637)),
    (@(45620.99999999999, 26, -1.65183899879832, 56.44606190205304)),
    (@(45627.99999999999, 26, -4.602482728117749, 55.06338686899744)),
    (@(45634.99999999999, 26, -3.897581428986077, 55.21639395440823)),
    (@(45641.99999999999, 26, -3.357611898020767, 55.04546018361148)),
    (@(45648.99999999999, 26, -5.488555514834304, 56.44835017869499)),
    (@(45655.99999999999, 26, -2.475925169480025, 55.29481478023038))
)

$data = New-Object 'object[,]' $rows.Count,4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $data[$i,0] = $rows[$i][0]
    $data[$i,1] = $rows[$i][1]
    $data[$i,2] = $rows[$i][2]
    $data[$i,3] = $rows[$i][3]
}

$ws3.Range("A2:D28").Value = $data

# Restore the originally active sheet/selection.
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null

Write-Output "PO Forecast sheet added."
